# Updated TPM-derived LR-pair values for Tnfsf8-Tnfrsf8 (rows 2-7)
# Row layout: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# followed by the 16 numeric NATMI metric columns (E..T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @{}
$data[2] = @("ECs", "Tnfsf8", "Tnfrsf8", "ECs", 1, 0.3333333333333333, 0.009345666666666667, 0.028037, 0.009908940382722529, 0.009908940382722529, 3, 1, 0.5079733333333333, 1.52392, 0.6037488273029235, 0.6037488273029235, 0.004747349448888888, 0.04272614504, 0.005982511135883308, 0.005982511135883308)
$data[3] = @("ECs", "Tnfsf8", "Tnfrsf8", "Resolving-Mac", 1, 0.3333333333333333, 0.009345666666666667, 0.028037, 0.009908940382722529, 0.009908940382722529, 2, 0.6666666666666666, 0.333392, 1.000176, 0.3962511726970765, 0.3962511726970764, 0.003115770501333333, 0.028041934512, 0.00392642924683922, 0.00392642924683922)
$data[4] = @("FAPs", "Tnfsf8", "Tnfrsf8", "ECs", 1, 0.3333333333333333, 0.1614316666666667, 0.484295, 0.1711613326194174, 0.1711613326194175, 3, 1, 0.5079733333333333, 1.52392, 0.6037488273029235, 0.6037488273029235, 0.08200298182222221, 0.7380268363999999, 0.1033384538485789, 0.1033384538485789)
$data[5] = @("FAPs", "Tnfsf8", "Tnfrsf8", "Resolving-Mac", 1, 0.3333333333333333, 0.1614316666666667, 0.484295, 0.1711613326194174, 0.1711613326194175, 2, 0.6666666666666666, 0.333392, 1.000176, 0.3962511726970765, 0.3962511726970764, 0.05382002621333333, 0.48438023592, 0.06782287877083854, 0.06782287877083854)
$data[6] = @("Resolving-Mac", "Tnfsf8", "Tnfrsf8", "ECs", 3, 1, 0.7723776666666667, 2.317133, 0.81892972699786, 0.81892972699786, 3, 1, 0.5079733333333333, 1.52392, 0.6037488273029235, 0.6037488273029235, 0.3923472579288889, 3.53112532136, 0.4944278623184613, 0.4944278623184613)
$data[7] = @("Resolving-Mac", "Tnfsf8", "Tnfrsf8", "Resolving-Mac", 3, 1, 0.7723776666666667, 2.317133, 0.81892972699786, 0.81892972699786, 2, 0.6666666666666666, 0.333392, 1.000176, 0.3962511726970765, 0.3962511726970764, 0.2575045350453333, 2.317540815408, 0.3245018646793987, 0.3245018646793987)

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $rowVals[$i]
    }
}